$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$rowCount = $usedRange.Rows.Count

for ($r = 1; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $val = $cell.Value2
    if ($val -ne $null -and $val -is [string] -and $val.EndsWith("_base")) {
        $cell.Value2 = $val.Substring(0, $val.Length - 5)
    }
}
